$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGDPbES")
$ws.Range("B11").Value = 0.106
$ws.Activate()
$ws.Range("B12").Select()
